# feat: add 2022-Q1 data
#
# 1. Duplicate the "2021-Q4" sheet (same column layout: 基金代码/基金名称/
#    基金规模/股票总仓位/仓位占比/持有市值(亿元)/仓位排名) to use as the new
#    "2022-Q1" sheet, inserted right before the "总计" summary sheet.
# 2. Overwrite its data with the 2022-Q1 fund holdings.
# 3. Insert a new top row into "总计" for the 2022-Q1 summary line.

$wb = $excel.ActiveWorkbook

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)          # "总计"
$q4 = $wb.Worksheets.Item("2021-Q4")

# Duplicate 2021-Q4 -> new sheet placed immediately before 总计, then rename.
$q4.Copy($lastSheet)

# Re-resolve "总计" by name (the old object reference no longer tracks the
# "总计" tab once the sheet collection has shifted) and grab its new
# left-hand neighbour - the freshly duplicated sheet.
$zjSheet = $wb.Worksheets.Item("总计")
$newSheet = $wb.Worksheets.Item($zjSheet.Index - 1)
$newSheet.Name = "2022-Q1"

# ---- fill the new 2022-Q1 detail sheet --------------------------------
# Row 1 header + column A index numbers already match (copied). Only the
# B:H data cells need new values.

function Set-TextCell($ws, $addr, $text) {
    # Force text (not number) interpretation for numeric-looking strings
    # (fund codes, decimal figures) the same way Excel does when a cell is
    # pre-formatted as Text, then drop the number-format override again so
    # the cell is left with the sheet's normal (unstyled) look.
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $text
    $ws.Range($addr).ClearFormats()
}

# Row 2 - 513090 易方达中证香港证券投资主题ETF
Set-TextCell $newSheet "B2" "513090"
Set-TextCell $newSheet "C2" "易方达中证香港证券投资主题ETF"
Set-TextCell $newSheet "D2" "11.07"
Set-TextCell $newSheet "E2" "96.47"
Set-TextCell $newSheet "F2" "6.29"
Set-TextCell $newSheet "G2" "0.6963"
$newSheet.Range("H2").Value = 7

# Row 3 - 004496 前海开源多元策略灵活配置混合A
Set-TextCell $newSheet "B3" "004496"
Set-TextCell $newSheet "C3" "前海开源多元策略灵活配置混合A"
Set-TextCell $newSheet "D3" "3.09"
Set-TextCell $newSheet "E3" "91.25"
Set-TextCell $newSheet "F3" "8.88"
Set-TextCell $newSheet "G3" "0.2744"
$newSheet.Range("H3").Value = 3

# Row 4 - 004497 前海开源多元策略灵活配置混合C
Set-TextCell $newSheet "B4" "004497"
Set-TextCell $newSheet "C4" "前海开源多元策略灵活配置混合C"
Set-TextCell $newSheet "D4" "2.06"
Set-TextCell $newSheet "E4" "91.25"
Set-TextCell $newSheet "F4" "8.88"
Set-TextCell $newSheet "G4" "0.1829"
$newSheet.Range("H4").Value = 3

$newSheet.Range("A1:H4").Select()

# ---- insert the 2022-Q1 summary row into 总计 --------------------------
$zj = $wb.Worksheets.Item("总计")
$zj.Rows(2).Insert()

$zj.Range("A2").Value = 0
$zj.Range("B2").Value = "2022-Q1"
$zj.Range("C2").Value = 3
$zj.Range("D2").Value = 1.15

$zj.Range("A2").Style = $zj.Range("A3").Style
